$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window / view tweaks (best effort; host may not persist these) ---
$win = $excel.Windows.Item(1)
$win.Left = -120
$win.Top = -120
$win.Width = 38640
$win.Height = 21120

# --- Column widths ---
# Column B (new "Member" column) gets its own width.
$ws.Columns.Item(2).ColumnWidth = 37.42578125
# Column J (AVERAGE column) becomes its own width, separate from C:I.
$ws.Columns.Item(10).ColumnWidth = 21.7109375

# --- New data rows 3-5 ---
$ws.Range("B3").Value = "Adrian Macauley s225733"
$ws.Range("C3").Value = 33
$ws.Range("D3").Value = 33
$ws.Range("E3").Value = 50
$ws.Range("J3").Formula = "=AVERAGE(C3:I3)"

$ws.Range("B4").Value = "Peter Juul s215781"
$ws.Range("C4").Value = 33
$ws.Range("D4").Value = 33
$ws.Range("E4").Value = 23

$ws.Range("B5").Value = "Sivalaxmanan B. Krishnapillai s245231"
$ws.Range("C5").Value = 33
$ws.Range("D5").Value = 33
$ws.Range("E5").Value = 45

# Shared formula across J4:J5
$ws.Range("J4:J5").Formula = "=AVERAGE(C4:I4)"

# --- Header cell J2: "AVERAGE" -> "AVERAGE in percentage" (renamed last so it
#     keeps the tail slot in sharedStrings, matching the authored order) ---
$ws.Range("J2").Value = "AVERAGE in percentage"

# --- Formatting for the last member's name cell ---
$ws.Range("B5").Font.Name = "Aptos"
$ws.Range("B5").VerticalAlignment = -4108

# --- Selection matches the saved workbook state ---
$ws.Range("L5").Select()
